$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set the "Who writes it" values for the Introduction/Methodology sub-rows (E6, E7) to "Franz"
$ws.Range("E6").Value = "Franz"
$ws.Range("E7").Value = "Franz"

# Update the active selection to E7 (matches the diff's selection change)
$ws.Range("E7").Select()
